$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New row 15 mirrors row 14's formatting (copy style from D14/E14 downwards)
$ws.Range("D14:E14").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D15").Value = 44266
$ws.Range("E15").Value = "Fin de la première version"

$ws.Range("T13").Select()
